$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.065.48"
$ws.Range("E2").Value = "  -5.04%  "

$ws.Range("D3").Value = "2.921.39"
$ws.Range("E3").Value = "  -2.57%  "

$ws.Range("E4").Value = "  +0.20%  "

$ws.Range("D5").Value = "553.68"
$ws.Range("E5").Value = "  -1.95%  "

$ws.Range("D6").Value = "124.09"
$ws.Range("E6").Value = "  -2.92%  "

$ws.Range("E7").Value = "  +0.09%  "

$ws.Range("D8").Value = "2.916.07"
$ws.Range("E8").Value = "  -2.54%  "

$ws.Range("D9").Value = "0.495"
$ws.Range("E9").Value = "  +0.44%  "

$ws.Range("D10").Value = "0.126"
$ws.Range("E10").Value = "  -5.68%  "

$ws.Range("D11").Value = "4.76"
$ws.Range("E11").Value = "  -7.09%  "

$ws.Range("D12").Value = "0.439"
$ws.Range("E12").Value = "  +2.08%  "

$ws.Range("D13").Value = "0.0000214"
$ws.Range("E13").Value = "  -4.29%  "

$ws.Range("D14").Value = "32.06"
$ws.Range("E14").Value = "  -2.34%  "

$ws.Range("E15").Value = "  +1.03%  "

$ws.Range("D16").Value = "3.403.73"
$ws.Range("E16").Value = "  -2.19%  "

$ws.Range("D17").Value = "2.912.22"
$ws.Range("E17").Value = "  -2.62%  "

$ws.Range("B18").Value = "Polkadot"
$ws.Range("C18").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D18").Value = "6.56"
$ws.Range("E18").Value = "  +5.58%  "

$ws.Range("B19").Value = "WrappedBTC"
$ws.Range("C19").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D19").Value = "58.015.16"
$ws.Range("E19").Value = "  -4.97%  "

$ws.Range("D20").Value = "410.72"
$ws.Range("E20").Value = "  -5.54%  "

$ws.Range("D21").Value = "13.00"
$ws.Range("E21").Value = "  -1.03%  "

$ws.Range("D22").Value = "0.664"
$ws.Range("E22").Value = "  +0.59%  "

$ws.Range("D23").Value = "6.86"
$ws.Range("E23").Value = "  -3.53%  "

$ws.Range("D24").Value = "12.80"
$ws.Range("E24").Value = "  +1.99%  "

$ws.Range("D25").Value = "77.43"
$ws.Range("E25").Value = "  -1.84%  "

$ws.Range("E26").Value = "  +0.21%  "

$ws.Range("E27").Value = "  +0.09%  "

$ws.Range("D28").Value = "2.48"
$ws.Range("E28").Value = "  -0.13%  "

$ws.Range("D29").Value = "7.31"
$ws.Range("E29").Value = "  +0.57%  "

$ws.Range("D30").Value = "1.94"
$ws.Range("E30").Value = "  +0.16%  "

$ws.Range("D31").Value = "6.11"
$ws.Range("E31").Value = "  -1.91%  "

$ws.Range("D32").Value = "24.84"
$ws.Range("E32").Value = "  -2.07%  "

$ws.Range("D33").Value = "0.0964"
$ws.Range("E33").Value = "  +2.93%  "

$ws.Range("D34").Value = "0.918"
$ws.Range("E34").Value = "  -3.36%  "

$ws.Range("D35").Value = "2.05"
$ws.Range("E35").Value = "  -10.08%  "

$ws.Range("D36").Value = "5.39"
$ws.Range("E36").Value = "  -2.76%  "

$ws.Range("D37").Value = "48.09"
$ws.Range("E37").Value = "  -4.02%  "

$ws.Range("D38").Value = "8.48"
$ws.Range("E38").Value = "  +9.88%  "

$ws.Range("D39").Value = "0.0₃0635"
$ws.Range("E39").Value = "  -6.28%  "

$ws.Range("D40").Value = "0.0349"
$ws.Range("E40").Value = "  -3.62%  "

$ws.Range("D41").Value = "0.107"
$ws.Range("E41").Value = "  -1.28%  "

$ws.Range("D42").Value = "2.626.58"
$ws.Range("E42").Value = "  -0.35%  "

$ws.Range("B43").Value = "dogwifhat"
$ws.Range("C43").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D43").Value = "2.43"
$ws.Range("E43").Value = "  -0.63%  "

$ws.Range("B44").Value = "Bittensor"
$ws.Range("C44").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D44").Value = "363.31"
$ws.Range("E44").Value = "  -2.04%  "

$ws.Range("E45").Value = "  +0.00%  "

$ws.Range("D46").Value = "0.231"
$ws.Range("E46").Value = "  -1.68%  "

$ws.Range("D47").Value = "118.25"
$ws.Range("E47").Value = "  -1.80%  "

$ws.Range("B48").Value = "Fetch.AI"
$ws.Range("C48").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D48").Value = "1.97"
$ws.Range("E48").Value = "  +0.40%  "

$ws.Range("B49").Value = "Stellar"
$ws.Range("C49").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D49").Value = "0.107"
$ws.Range("E49").Value = "  +0.73%  "

$ws.Range("D50").Value = "22.98"
$ws.Range("E50").Value = "  -1.60%  "

$ws.Range("D51").Value = "1.97"
$ws.Range("E51").Value = "  -2.20%  "
